# Ajout d'un pdf de test avec conjecture et formule mathematique
#
# - Articles!A24:H24  -> new row (Id 22, "Test" article, test_overleaf.pdf)
# - Conjectures!A24:C24 -> matching "aucune conjecture (json manquant)" row
# - Column widths on the touched columns get (re)fit, same as Excel does
#   automatically when new, wider/narrower content lands in a column.
# - Selections end up on the last-touched cell of each sheet, with
#   "Conjectures" left as the active (visible) tab, matching the workbook's
#   activeTab=1.

$wb = $excel.ActiveWorkbook

$wsArticles    = $wb.Worksheets.Item("Articles")
$wsConjectures = $wb.Worksheets.Item("Conjectures")

# --- Articles: new row 24 -------------------------------------------------
$wsArticles.Cells.Item(24, 1).Value = 22
$wsArticles.Cells.Item(24, 2).Value = "Test"
$wsArticles.Cells.Item(24, 3).Value = "Sami"
$wsArticles.Cells.Item(24, 4).Value = "AA"
$wsArticles.Cells.Item(24, 5).Value = "AA"
$wsArticles.Cells.Item(24, 6).Value = "AA"
$wsArticles.Cells.Item(24, 7).Value = "AA"
$wsArticles.Cells.Item(24, 8).Value = "test_overleaf.pdf"

# --- Conjectures: matching row 24 ------------------------------------------
$wsConjectures.Cells.Item(24, 1).Value = 22
$wsConjectures.Cells.Item(24, 2).Value = "aucune conjecture (json manquant)"
$wsConjectures.Cells.Item(24, 3).Value = ""

# --- Column widths: re-fit the columns whose longest entry changed --------
$wsArticles.Columns.Item(2).ColumnWidth = 112.42578125
$wsArticles.Columns.Item(7).ColumnWidth = 31.5703125
$wsArticles.Columns.Item(8).ColumnWidth = 135.7109375

$wsConjectures.Columns.Item(1).ColumnWidth = 9.5703125
$wsConjectures.Columns.Item(2).ColumnWidth = 53.85546875

# --- Selections -------------------------------------------------------------
# Touch Articles last cell first, then Conjectures, so Conjectures stays the
# active/visible tab (matching the saved workbook's activeTab).
$wsArticles.Activate()
$wsArticles.Range("G24").Select() | Out-Null

$wsConjectures.Activate()
$wsConjectures.Range("A22").Select() | Out-Null

Write-Output "Added test article row (Articles!A24:H24) and matching Conjectures!A24:C24."
